$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$cells = @("I2","J2","L2","I3","J3","I4","J4","I5","J5","I6","J6","I7","J7","I8","J8","I9","J9","I10","J10","I11","J11","I12","J12")
foreach ($c in $cells) {
    $ws.Range($c).ClearContents()
}

$ws.Range("M8").Select()
